# Apply updated cryptocurrency price/volume data to Sheet1
# Commit: "Updated cryptos list on Sun Jun 23 16:43:42 UTC 2024 with GitHub Actions"
#
# Some Price (column D) values are purely numeric-looking strings (e.g. "584.65").
# Excel's COM layer would otherwise auto-convert these to floating point numbers
# (introducing binary-float rounding noise and losing formatting such as trailing
# zeros), so for those specific cells we force the cell's number format to Text
# ("@") before assigning the value, guaranteeing the value round-trips as the
# exact original text. Cells whose new text is not a bare number (e.g. multi-dot
# "64.120.40", or the whitespace-padded percentage strings in column E) do not
# need this treatment since Excel already keeps them as text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '64.120.40'
$ws.Range('E2').Value = '  -0.26%  '
$ws.Range('D3').Value = '3.479.21'
$ws.Range('E3').Value = '  -0.58%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '584.65'
$ws.Range('E5').Value = '  -0.22%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '131.62'
$ws.Range('E6').Value = '  -1.99%  '
$ws.Range('E7').Value = '  +0.16%  '
$ws.Range('E8').Value = '  -0.96%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '7.63'
$ws.Range('E9').Value = '  +4.99%  '
$ws.Range('E10').Value = '  -1.79%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.386'
$ws.Range('E11').Value = '  -0.54%  '
$ws.Range('D12').Value = '4.070.16'
$ws.Range('E12').Value = '  -0.61%  '
$ws.Range('E14').Value = '  -2.67%  '
$ws.Range('D15').Value = '3.480.30'
$ws.Range('E15').Value = '  -0.49%  '
$ws.Range('D16').Value = '64.124.66'
$ws.Range('E16').Value = '  -0.27%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '24.34'
$ws.Range('E17').Value = '  -6.63%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '9.97'
$ws.Range('E18').Value = '  +0.31%  '
$ws.Range('E19').Value = '  -1.26%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '13.43'
$ws.Range('E20').Value = '  -2.18%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '384.66'
$ws.Range('E21').Value = '  -2.51%  '
$ws.Range('E22').Value = '  -0.35%  '
$ws.Range('D23').Value = '3.618.51'
$ws.Range('E23').Value = '  -0.58%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '74.78'
$ws.Range('E24').Value = '  +0.78%  '
$ws.Range('E25').Value = '  +0.12%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '5.63'
$ws.Range('E26').Value = '  -0.53%  '
$ws.Range('E27').Value = '  -2.15%  '
$ws.Range('E28').Value = '  +0.04%  '
$ws.Range('E29').Value = '  -0.04%  '
$ws.Range('E30').Value = '  -4.17%  '
$ws.Range('E31').Value = '  -4.28%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '7.94'
$ws.Range('E32').Value = '  -4.60%  '
$ws.Range('D33').Value = '3.508.09'
$ws.Range('E33').Value = '  -0.28%  '
$ws.Range('E34').Value = '  +0.85%  '
$ws.Range('E35').Value = '  -0.01%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '22.96'
$ws.Range('E36').Value = '  -2.45%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.20'
$ws.Range('E37').Value = '  -0.07%  '
$ws.Range('E38').Value = '  -2.42%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.50'
$ws.Range('E39').Value = '  -4.12%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '162.65'
$ws.Range('E40').Value = '  +0.56%  '
$ws.Range('E41').Value = '  -1.14%  '
$ws.Range('E42').Value = '  -1.21%  '
$ws.Range('E43').Value = '  -0.06%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '41.38'
$ws.Range('E44').Value = '  -1.06%  '
$ws.Range('E45').Value = '  -2.78%  '
$ws.Range('E46').Value = '  -1.99%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '23.38'
$ws.Range('E47').Value = '  -6.96%  '
$ws.Range('E48').Value = '  -3.57%  '
$ws.Range('E49').Value = '  -1.35%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.904'
$ws.Range('E50').Value = '  +0.69%  '
$ws.Range('D51').Value = '2.350.62'
$ws.Range('E51').Value = '  -4.77%  '
